$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 21
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 10
$ws.Range("AB2").Value = 19
$ws.Range("AH2").Value = 26
$ws.Range("AK2").Value = 81
$ws.Range("AX2").Value = 34

# Row 3
$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 4.5
$ws.Range("K3").Value = 2.4
$ws.Range("L3").Value = 4.5
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 15
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 4.5
$ws.Range("Q3").Value = 1.62
$ws.Range("R3").Value = 2.25
$ws.Range("S3").Value = 1.29
$ws.Range("T3").Value = 3.5
$ws.Range("U3").Value = 1.57
$ws.Range("V3").Value = 2.25
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 9.5
$ws.Range("AA3").Value = 13
$ws.Range("AB3").Value = 21
$ws.Range("AC3").Value = 15
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 41
$ws.Range("AG3").Value = 126
$ws.Range("AH3").Value = 17
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 34
$ws.Range("AP3").Value = 17
$ws.Range("AS3").Value = 101
$ws.Range("AT3").Value = 3.5
$ws.Range("AU3").Value = 7.5
$ws.Range("AY3").Value = 26
$ws.Range("BA3").Value = 81
$ws.Range("BB3").Value = 151
$ws.Range("BC3").Value = 351

# Row 4
$ws.Range("G4").Value = 2.1
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.75
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 9.5
$ws.Range("AC4").Value = 9
$ws.Range("AE4").Value = 15
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 19
$ws.Range("AM4").Value = 41
$ws.Range("AP4").Value = 23
$ws.Range("AT4").Value = 2.63
$ws.Range("AX4").Value = 21
$ws.Range("BA4").Value = 101

# Row 12
$ws.Range("M12").Value = 1.03
$ws.Range("N12").Value = 10

# Row 15
$ws.Range("G15").Value = 2.1
$ws.Range("I15").Value = 3.75
$ws.Range("K15").Value = 2.05
$ws.Range("L15").Value = 4.33
$ws.Range("N15").Value = 7.5
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 1.65
$ws.Range("X15").Value = 9.5
$ws.Range("Z15").Value = 19
$ws.Range("AG15").Value = 301
$ws.Range("AH15").Value = 10
$ws.Range("AI15").Value = 19
$ws.Range("AL15").Value = 34
$ws.Range("AO15").Value = 12
$ws.Range("AP15").Value = 23
$ws.Range("AW15").Value = 5.5
$ws.Range("AY15").Value = 34

# Row 16
$ws.Range("G16").Value = 2.1
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 3.8
$ws.Range("J16").Value = 2.88
$ws.Range("K16").Value = 1.91
$ws.Range("L16").Value = 4.75
$ws.Range("M16").Value = 1.13
$ws.Range("N16").Value = 6
$ws.Range("O16").Value = 1.5
$ws.Range("P16").Value = 2.5
$ws.Range("Q16").Value = 2.6
$ws.Range("R16").Value = 1.48
$ws.Range("S16").Value = 1.57
$ws.Range("T16").Value = 2.25
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 8.5
$ws.Range("Y16").Value = 10
$ws.Range("Z16").Value = 19
$ws.Range("AA16").Value = 21
$ws.Range("AH16").Value = 8.5
$ws.Range("AI16").Value = 17
$ws.Range("AJ16").Value = 15
$ws.Range("AL16").Value = 41
$ws.Range("AN16").Value = 4
$ws.Range("AO16").Value = 13
$ws.Range("AP16").Value = 29
$ws.Range("AQ16").Value = 41
$ws.Range("AR16").Value = 81
$ws.Range("AT16").Value = 2.25
$ws.Range("AW16").Value = 5.5
$ws.Range("AX16").Value = 23
$ws.Range("AY16").Value = 41
$ws.Range("AZ16").Value = 81

Write-Output "Applied all odds updates"